$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 713688901027.9006
$ws.Range("C3").Value = 240904302203.1317
$ws.Range("C4").Value = 37421036288.88834
$ws.Range("C5").Value = 33516769668.61184
$ws.Range("C6").Value = 25941175676.57832
$ws.Range("C7").Value = 13491724901.42055
$ws.Range("C8").Value = 10948878600.2546
$ws.Range("C9").Value = 9133398904.990887
$ws.Range("C10").Value = 8162374780.375494
$ws.Range("C11").Value = 8110987920.436409
$ws.Range("C12").Value = 7935764574.283823
$ws.Range("C13").Value = 7888159434.931012
$ws.Range("C14").Value = 7004567269.532668
$ws.Range("C15").Value = 5970710377.514452
$ws.Range("C16").Value = 5305017541.105283
$ws.Range("C17").Value = 5127623626.674
$ws.Range("C18").Value = 4542817292.301433
$ws.Range("C19").Value = 3760213762.134553
$ws.Range("C20").Value = 3664101041.967874
$ws.Range("C21").Value = 3485833767.538836
$ws.Range("C22").Value = 3351316873.553715
$ws.Range("C23").Value = 2964508064.322744
$ws.Range("C24").Value = 2928518666.515219
$ws.Range("C25").Value = 2835973127.381089
$ws.Range("C26").Value = 2482378317.382869
